$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# B2 ("Bitcoin") carries the workbook default (unstyled) format; used to
# restore that default style on D-column cells after forcing them to text
# below (quote-prefixing a numeric-looking value stamps a quotePrefix style).
$origStyle = $ws.Range("B2").Style

$ws.Range('D2').Value = '60.227.43'
$ws.Range('E2').Value = '  +3.81%  '
$ws.Range('D3').Value = '3.199.21'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'538.35"
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('D6').Value = "'145.54"
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +4.86%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = "'0.519"
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  +2.90%  '
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('E10').Value = '  +4.79%  '
$ws.Range('E11').Value = '  +2.80%  '
$ws.Range('D12').Value = '3.749.74'
$ws.Range('E12').Value = '  +2.16%  '
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = "'0.0000175"
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +3.34%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = "'26.13"
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +1.78%  '
$ws.Range('D16').Value = '60.227.12'
$ws.Range('E16').Value = '  +3.67%  '
$ws.Range('D17').Value = '3.210.26'
$ws.Range('E17').Value = '  +2.12%  '
$ws.Range('D18').Value = "'6.20"
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('E19').Value = '  +1.48%  '
$ws.Range('D20').Value = "'8.40"
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +3.19%  '
$ws.Range('D21').Value = "'383.87"
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +2.40%  '
$ws.Range('D22').Value = "'1.00"
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('E23').Value = '  +3.44%  '
$ws.Range('D24').Value = "'70.10"
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = "'0.172"
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  +2.76%  '
$ws.Range('D26').Value = "'8.84"
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +13.47%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').Value = '0.0₃0904'
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('E29').Value = '  +1.82%  '
$ws.Range('D30').Value = "'22.49"
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +3.38%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').Value = '  +5.52%  '
$ws.Range('D33').Value = "'1.22"
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +3.43%  '
$ws.Range('D34').Value = "'6.54"
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +5.19%  '
$ws.Range('D35').Value = "'156.73"
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  -2.57%  '
$ws.Range('D36').Value = "'1.36"
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '2.775.66'
$ws.Range('E37').Value = '  +7.10%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = "'25.81"
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +1.46%  '
$ws.Range('D39').Value = "'0.0716"
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +5.92%  '
$ws.Range('D40').Value = "'1.71"
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +1.35%  '
$ws.Range('D41').Value = "'4.27"
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('D42').Value = "'39.76"
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +2.52%  '
$ws.Range('D43').Value = "'0.730"
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  +4.31%  '
$ws.Range('E44').Value = '  +6.13%  '
$ws.Range('E45').Value = '  +3.39%  '
$ws.Range('D46').Value = '3.241.37'
$ws.Range('E46').Value = '  +2.16%  '
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('E48').Value = '  +2.19%  '
$ws.Range('D49').Value = "'0.802"
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +7.29%  '
$ws.Range('D50').Value = "'20.61"
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +2.41%  '
$ws.Range('D51').Value = "'1.00"
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +0.00%  '
